$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume text values (which often look numeric, e.g. "0.9991")
# are stored as text, matching the original inlineStr cells, not auto-coerced
# to numbers by Excel's normal type inference.
$textRange = $ws.Range("B2:E51")
$textRange.NumberFormat = "@"

$ws.Range("D2").Value = '29.379.97'
$ws.Range("E2").Value = '  -0.01%  '
$ws.Range("D3").Value = '1.848.74'
$ws.Range("E3").Value = '  +0.02%  '
$ws.Range("D4").Value = '0.9991'
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '240.27'
$ws.Range("E5").Value = '  -0.07%  '
$ws.Range("D6").Value = '0.6284'
$ws.Range("E6").Value = '  -0.16%  '
$ws.Range("D7").Value = '0.9997'
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").Value = '0.07626'
$ws.Range("E8").Value = '  -0.21%  '
$ws.Range("E9").Value = '  -1.13%  '
$ws.Range("D10").Value = '24.75'
$ws.Range("E10").Value = '  +1.03%  '
$ws.Range("D11").Value = '0.07735'
$ws.Range("E11").Value = '  -0.16%  '
$ws.Range("D12").Value = '5.032'
$ws.Range("E12").Value = '  +0.50%  '
$ws.Range("D13").Value = '0.6796'
$ws.Range("E14").Value = '  -3.15%  '
$ws.Range("E15").Value = '  -0.40%  '
$ws.Range("E16").Value = '  +0.58%  '
$ws.Range("D17").Value = '29.396.13'
$ws.Range("E17").Value = '  -0.08%  '
$ws.Range("D18").Value = '228.02'
$ws.Range("E18").Value = '  -0.50%  '
$ws.Range("E19").Value = '  -0.73%  '
$ws.Range("D20").Value = '0.9994'
$ws.Range("E20").Value = '  -0.08%  '
$ws.Range("D21").Value = '7.486'
$ws.Range("E21").Value = '  +0.64%  '
$ws.Range("D22").Value = '0.9999'
$ws.Range("E22").Value = '  -0.07%  '
$ws.Range("D23").Value = '158.84'
$ws.Range("D24").Value = '0.1387'
$ws.Range("E24").Value = '  -0.18%  '
$ws.Range("D25").Value = '8.406'
$ws.Range("E25").Value = '  +0.32%  '
$ws.Range("D26").Value = '17.69'
$ws.Range("E26").Value = '  +0.18%  '
$ws.Range("D27").Value = '1.407'
$ws.Range("E27").Value = '  +8.79%  '
$ws.Range("E28").Value = '  -0.23%  '
$ws.Range("D29").Value = '0.05599'
$ws.Range("E29").Value = '  -0.65%  '
$ws.Range("E30").Value = '  +0.15%  '
$ws.Range("D31").Value = '4.070'
$ws.Range("E31").Value = '  +0.63%  '
$ws.Range("E32").Value = '  +0.62%  '
$ws.Range("D33").Value = '1.835'
$ws.Range("E33").Value = '  -0.77%  '
$ws.Range("E34").Value = '  -1.23%  '
$ws.Range("D35").Value = '2.582'
$ws.Range("E35").Value = '  -0.15%  '
$ws.Range("B36").Value = 'VeChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D36").Value = '0.01807'
$ws.Range("B37").Value = 'Maker'
$ws.Range("C37").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D37").Value = '1.232.68'
$ws.Range("E37").Value = '  +0.27%  '
$ws.Range("D38").Value = '2.718'
$ws.Range("E38").Value = '  -2.00%  '
$ws.Range("D39").Value = '6.382'
$ws.Range("E39").Value = '  -1.01%  '
$ws.Range("D40").Value = '0.9025'
$ws.Range("E40").Value = '  -0.76%  '
$ws.Range("D41").Value = '0.9998'
$ws.Range("E41").Value = '  -0.01%  '
$ws.Range("D42").Value = '101.50'
$ws.Range("E42").Value = '  +0.09%  '
$ws.Range("D43").Value = '66.04'
$ws.Range("E43").Value = '  -0.05%  '
$ws.Range("D44").Value = '7.203'
$ws.Range("E44").Value = '  +0.70%  '
$ws.Range("B45").Value = 'BabyDogeCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D45").Value = '0.00000000118'
$ws.Range("E45").Value = '  -2.82%  '
$ws.Range("B46").Value = 'TheSandbox'
$ws.Range("C46").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D46").Value = '0.4002'
$ws.Range("E46").Value = '  -0.14%  '
$ws.Range("D47").Value = '9.011'
$ws.Range("E47").Value = '  -0.37%  '
$ws.Range("D48").Value = '1.677'
$ws.Range("E48").Value = '  -0.56%  '
$ws.Range("D49").Value = '0.1135'
$ws.Range("E49").Value = '  +1.01%  '
$ws.Range("D50").Value = '0.05705'
$ws.Range("E51").Value = '  +0.05%  '

# Remove the temporary Text number-format so cell styling matches the original
# (unstyled) cells; the values already committed as text remain text.
$textRange.ClearFormats()
